# List View - All implemented
# Adds a new "Listviews are cropped at the top and bottom" row (row 18)
# to the Documentation sheet, describing the fix with top/bottom padding
# notes, and normalizes the wrap-text styling on column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply wrap text to the existing D16 cell so its style slot collapses
# onto the shared wrap-text style already used elsewhere in column D/C.
$ws.Range("D16").WrapText = $true

# New row 18: Feature "A" - cropped listviews issue + fix notes.
$ws.Range("B18").Value = "A"
$ws.Range("C18").Value = "Listviews are cropped at the top and bottom"
$ws.Range("D18").Value = "add a top padding to the listviews of taskQ_dialog_padding (16dp)`nadd a bottom padding to the listviews of tapp_icon_size (48dp), which must be same as hight of tab_Button_NewItem button`n"

$ws.Range("C18").WrapText = $true
$ws.Range("D18").WrapText = $true

# Move the active selection the way it ended up after the edit.
$null = $ws.Range("D24").Select()
